# Generate Report for Archive
#
# The nightly localization-status report is regenerated: every row that was
# previously "Ready for handoff" has moved on to "In Translation", and after
# the (shorter) new text is written, Excel re-autofits the status columns so
# they hug their content again.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- "Overview" sheet: zh-cn / de-de status columns (E and F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$usedRange = $wsOverview.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    if ($wsOverview.Cells.Item($r, 5).Value2 -eq $oldStatus) {
        $wsOverview.Cells.Item($r, 5).Value = $newStatus
    }
    if ($wsOverview.Cells.Item($r, 6).Value2 -eq $oldStatus) {
        $wsOverview.Cells.Item($r, 6).Value = $newStatus
    }
}
$wsOverview.Columns("E:F").AutoFit()
$wsOverview.Columns("E:F").ColumnWidth = 12.5

# --- per-language sheets: "Status" column (C) ---
$langSheets = @("zh-cn", "de-de")
foreach ($sheetName in $langSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    $last = $used.Rows.Count

    for ($r = 2; $r -le $last; $r++) {
        if ($ws.Cells.Item($r, 3).Value2 -eq $oldStatus) {
            $ws.Cells.Item($r, 3).Value = $newStatus
        }
    }
    $ws.Columns("C").AutoFit()
    $ws.Columns("C").ColumnWidth = 12.5
}
